$wb = $excel.ActiveWorkbook

# --- Filter sheet: fix the wording of the 1987/1990-1991 filter-step label ---
# (comma was outside the parenthetical; move it inside, matching the other labels)
$filterWs = $wb.Worksheets.Item("Filter")
$filterWs.Range("A3").Value = "Remove 1987 (sampled only MNT), 1990-1991 low sample sizes"

# Give the Drifts / PositiveDrifts count columns an integer number format
$filterWs.Range("B1:C7").NumberFormat = "0"

# Match column widths for the two numeric columns
$filterWs.Columns("B:C").ColumnWidth = 8.25

# Switch the sheet to portrait for printing
$filterWs.PageSetup.Orientation = 1

# Make "Filter" the active sheet/tab, with A7 selected (matches the saved view)
$filterWs.Activate() | Out-Null
$filterWs.Range("A7").Select() | Out-Null
